# Auto-generated edit script: updates cryptos D (Price) and E (Volume(1h)) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.091.26"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "1.653.72"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").Value = "'218.70"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").Value = "'0.5247"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("D8").Value = "'0.2681"
$ws.Range("E8").Value = "  +1.33%  "
$ws.Range("D9").Value = "'0.06360"
$ws.Range("D10").Value = "'20.50"
$ws.Range("E10").Value = "  -1.83%  "
$ws.Range("D11").Value = "'0.07677"
$ws.Range("E11").Value = "  -1.90%  "
$ws.Range("D12").Value = "'4.595"
$ws.Range("E12").Value = "  +1.64%  "
$ws.Range("D13").Value = "1.643.41"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("D14").Value = "1.881.43"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").Value = "'0.5611"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "0.0₅8243"
$ws.Range("E16").Value = "  +1.77%  "
$ws.Range("D17").Value = "'65.63"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "26.081.70"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").Value = "'4.683"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").Value = "'10.32"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").Value = "'189.68"
$ws.Range("E22").Value = "  -5.15%  "
$ws.Range("D23").Value = "'5.973"
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("E24").Value = "  -0.60%  "
$ws.Range("D25").Value = "'146.10"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").Value = "'0.1199"
$ws.Range("E26").Value = "  -1.10%  "
$ws.Range("D27").Value = "'7.241"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "'15.93"
$ws.Range("E28").Value = "  -1.68%  "
$ws.Range("D29").Value = "'1.524"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("D30").Value = "'0.05629"
$ws.Range("E30").Value = "  -4.53%  "
$ws.Range("D31").Value = "'1.271"
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").Value = "'3.488"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("E33").Value = "  +1.54%  "
$ws.Range("D34").Value = "'1.580"
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("D37").Value = "'2.410"
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("D38").Value = "'0.5751"
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("D39").Value = "'0.01588"
$ws.Range("E39").Value = "  -1.55%  "
$ws.Range("D40").Value = "'5.967"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("D42").Value = "'0.8384"
$ws.Range("E42").Value = "  -2.27%  "
$ws.Range("D43").Value = "1.023.47"
$ws.Range("E43").Value = "  -4.80%  "
$ws.Range("D44").Value = "'101.22"
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").Value = "1.791.67"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("D46").Value = "'58.21"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("D47").Value = "0.0₈105"
$ws.Range("E47").Value = "  +3.42%  "
$ws.Range("E48").Value = "  -1.06%  "
$ws.Range("D49").Value = "'0.05304"
$ws.Range("E49").Value = "  +3.05%  "
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("D51").Value = "'0.4342"
$ws.Range("E51").Value = "  -1.72%  "
